$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("A144").Value = [double]"45930.49350694445"
$ws.Range("A144").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B144").Value = "0x01,0x7c"
$ws.Range("C144").Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Range("D144").Value = "0x00,0xC8"
$ws.Range("E144").Value = "0xf"
$ws.Range("F144").Value = 380
$ws.Range("G144").Value = [double]"7.598631275147109e+23"
$ws.Range("H144").Value = 208
$ws.Range("I144").Value = 15
$ws.Range("A145").Value = [double]"45931.49167824074"
$ws.Range("A145").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B145").Value = "0x01,0x7c"
$ws.Range("C145").Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Range("D145").Value = "0x00,0xC8"
$ws.Range("E145").Value = "0xf"
$ws.Range("F145").Value = 380
$ws.Range("G145").Value = [double]"7.598631275147109e+23"
$ws.Range("H145").Value = 208
$ws.Range("I145").Value = 15
$ws.Range("A146").Value = [double]"45932.49309027778"
$ws.Range("A146").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B146").Value = "0x01,0x7c"
$ws.Range("C146").Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Range("D146").Value = "0x00,0xC8"
$ws.Range("E146").Value = "0xf"
$ws.Range("F146").Value = 380
$ws.Range("G146").Value = [double]"7.598631275147109e+23"
$ws.Range("H146").Value = 208
$ws.Range("I146").Value = 15

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("A144").Value = [double]"45930.49350694445"
$ws.Range("A144").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B144").Value = "0x01,0x90"
$ws.Range("C144").Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Range("D144").Value = "0x00,0xCC"
$ws.Range("E144").Value = "0xe"
$ws.Range("F144").Value = 400
$ws.Range("G144").Value = [double]"5.68432987514711e+23"
$ws.Range("H144").Value = 216
$ws.Range("I144").Value = 14
$ws.Range("A145").Value = [double]"45931.49167824074"
$ws.Range("A145").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B145").Value = "0x01,0x90"
$ws.Range("C145").Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Range("D145").Value = "0x00,0xCC"
$ws.Range("E145").Value = "0xe"
$ws.Range("F145").Value = 400
$ws.Range("G145").Value = [double]"5.68432987514711e+23"
$ws.Range("H145").Value = 212
$ws.Range("I145").Value = 14
$ws.Range("A146").Value = [double]"45932.49309027778"
$ws.Range("A146").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B146").Value = "0x01,0x90"
$ws.Range("C146").Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Range("D146").Value = "0x00,0xCC"
$ws.Range("E146").Value = "0xe"
$ws.Range("F146").Value = 400
$ws.Range("G146").Value = [double]"5.68432987514711e+23"
$ws.Range("H146").Value = 212
$ws.Range("I146").Value = 14

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("A144").Value = [double]"45930.49350694445"
$ws.Range("A144").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B144").Value = "0x00,0x6e"
$ws.Range("C144").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Range("D144").Value = "0x00,0x53"
$ws.Range("E144").Value = "0x3"
$ws.Range("F144").Value = 110
$ws.Range("G144").Value = [double]"5.68631262647114e+23"
$ws.Range("H144").Value = 83
$ws.Range("I144").Value = 3
$ws.Range("A145").Value = [double]"45931.49167824074"
$ws.Range("A145").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B145").Value = "0x00,0x6e"
$ws.Range("C145").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Range("D145").Value = "0x00,0x53"
$ws.Range("E145").Value = "0x3"
$ws.Range("F145").Value = 110
$ws.Range("G145").Value = [double]"5.68631262647114e+23"
$ws.Range("H145").Value = 83
$ws.Range("I145").Value = 3
$ws.Range("A146").Value = [double]"45932.49309027778"
$ws.Range("A146").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B146").Value = "0x00,0x6e"
$ws.Range("C146").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Range("D146").Value = "0x00,0x53"
$ws.Range("E146").Value = "0x3"
$ws.Range("F146").Value = 110
$ws.Range("G146").Value = [double]"5.68631262647114e+23"
$ws.Range("H146").Value = 83
$ws.Range("I146").Value = 3

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("A144").Value = [double]"45930.49350694445"
$ws.Range("A144").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B144").Value = "0x00,0x6e"
$ws.Range("C144").Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Range("D144").Value = "0x00,0x50"
$ws.Range("E144").Value = "0x3"
$ws.Range("F144").Value = 110
$ws.Range("G144").Value = [double]"9.85046333984776e+23"
$ws.Range("H144").Value = 80
$ws.Range("I144").Value = 3
$ws.Range("A145").Value = [double]"45931.49167824074"
$ws.Range("A145").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B145").Value = "0x00,0x6e"
$ws.Range("C145").Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Range("D145").Value = "0x00,0x50"
$ws.Range("E145").Value = "0x3"
$ws.Range("F145").Value = 110
$ws.Range("G145").Value = [double]"9.85046333984776e+23"
$ws.Range("H145").Value = 80
$ws.Range("I145").Value = 3
$ws.Range("A146").Value = [double]"45932.49309027778"
$ws.Range("A146").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B146").Value = "0x00,0x6e"
$ws.Range("C146").Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Range("D146").Value = "0x00,0x50"
$ws.Range("E146").Value = "0x3"
$ws.Range("F146").Value = 110
$ws.Range("G146").Value = [double]"9.85046333984776e+23"
$ws.Range("H146").Value = 80
$ws.Range("I146").Value = 3

